$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64, pushing existing rows 64-125 down to 65-126.
$insertRow = $ws.Rows.Item(64)
$insertRow.Insert()

# Populate the newly inserted row 64 with the new data record.
$ws.Cells.Item(64, 1).Value = 1
$ws.Cells.Item(64, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(64, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(64, 4).Value = 45280
$ws.Cells.Item(64, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(64, 5).Value = 15
$ws.Cells.Item(64, 6).Value = 100112012
$ws.Cells.Item(64, 7).Value = "Espinaca"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 300
$ws.Cells.Item(64, 11).Value = 800
$ws.Cells.Item(64, 12).Value = 1000
$ws.Cells.Item(64, 13).Value = 900
$ws.Cells.Item(64, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(64, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(64, 16).Value = 300
$ws.Cells.Item(64, 17).Value = 3
$ws.Cells.Item(64, 18).Value = "Hortaliza"
